$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-34 (only cells whose content changed) ---
$ws.Range("B2").Value = "NSE:ABCAPITAL"
$ws.Range("D2").Value = "NSE:BHEL"
$ws.Range("E2").Value = "NSE:LTTS"
$ws.Range("F2").Value = "NSE:ABCAPITAL"

$ws.Range("B3").Value = "NSE:AEROFLEX"
$ws.Range("C3").Value = "NSE:DCMSHRIRAM"
$ws.Range("D3").Value = "NSE:IEX"
$ws.Range("F3").Value = "NSE:GMRINFRA"

$ws.Range("B4").Value = "NSE:AGARIND"
$ws.Range("C4").Value = "NSE:GILLETTE"
$ws.Range("D4").Value = "NSE:NMDC"
$ws.Range("F4").Value = "NSE:POLYCAB"

$ws.Range("B5").Value = "NSE:AROGRANITE"
$ws.Range("C5").Value = "NSE:HOMEFIRST"
$ws.Range("D5").Value = "NSE:RECLTD"

$ws.Range("B6").Value = "NSE:ASHAPURMIN"
$ws.Range("C6").Value = "NSE:JINDALPHOT"

$ws.Range("B7").Value = "NSE:ASPINWALL"
$ws.Range("C7").ClearContents()

$ws.Range("B8").Value = "NSE:AUTOIND"
$ws.Range("C8").ClearContents()

$ws.Range("B9").Value = "NSE:BANKINDIA"
$ws.Range("C9").ClearContents()

$ws.Range("B10").Value = "NSE:BLAL"
$ws.Range("C10").ClearContents()

$ws.Range("B11").Value = "NSE:BLISSGVS"
$ws.Range("C11").ClearContents()

$ws.Range("B12").Value = "NSE:DCMSRIND"

$ws.Range("B13").Value = "NSE:DEEPAKFERT"

$ws.Range("B14").Value = "NSE:DHRUV"

$ws.Range("B15").Value = "NSE:DMCC"

$ws.Range("B16").Value = "NSE:DPSCLTD"

$ws.Range("B17").Value = "NSE:EDELWEISS"

$ws.Range("B18").Value = "NSE:FINOPB"

$ws.Range("B19").Value = "NSE:GMRINFRA"

$ws.Range("B20").Value = "NSE:GNFC"

$ws.Range("B21").Value = "NSE:GRAPHITE"

$ws.Range("B22").Value = "NSE:HCC"

$ws.Range("B23").Value = "NSE:HEXATRADEX"

$ws.Range("B24").Value = "NSE:HINDCON"

$ws.Range("B25").Value = "NSE:HIRECT"

$ws.Range("B26").Value = "NSE:HONDAPOWER"

$ws.Range("B27").Value = "NSE:HPIL"

$ws.Range("B28").Value = "NSE:IRFC"

$ws.Range("B29").Value = "NSE:IVC"

$ws.Range("B30").Value = "NSE:JPPOWER"

$ws.Range("B31").Value = "NSE:KHADIM"

$ws.Range("B32").Value = "NSE:MAHABANK"

$ws.Range("B33").Value = "NSE:MANINFRA"

$ws.Range("B34").Value = "NSE:MOIL"

# --- Extend formatting (bold/border/centered style) for column A down to new rows ---
$ws.Range("A34").Copy() | Out-Null
$ws.Range("A35:A46").PasteSpecial(-4122) | Out-Null

# --- Add new rows 35-46 ---
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "NSE:MUKANDLTD"
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "NSE:PENIND"
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "NSE:PILITA"
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = "NSE:POLYCAB"
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = "NSE:PREMEXPLN"
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "NSE:RAJMET"
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "NSE:RCF"
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "NSE:ROML"
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = "NSE:RPGLIFE"
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "NSE:RPSGVENT"
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "NSE:RTNPOWER"
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = "NSE:SAKSOFT"

$excel.CutCopyMode = $false
